# RPA datasets push 2024-01-03
$wb = $excel.ActiveWorkbook

# --- Sheet 1 (01_IB전략컨설팅부): remove the last three IPO entries ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A26:Y28").EntireRow.Delete()

# --- Sheet 2 (02_38커뮤니케이션): insert the new 스튜디오삼익 demand-forecast
#     entry (updated schedule) ahead of IBKS스팩24호, and drop its old entry ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Rows.Item(4).Insert()
$ws2.Range("A4").Value = "스튜디오삼익"
$ws2.Range("B4").Value = "2024.01.17~01.23"
$ws2.Range("C4").Value = "14,500~16,500"
$ws2.Range("D4").Value = "-"
$ws2.Range("E4").Value = 12325
$ws2.Range("F4").Value = "DB금융투자"
$ws2.Rows.Item(14).Delete()
